$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (ALC) - hunk 0
$ws.Range("H33").Value = 448.36
$ws.Range("J33").Value = 913.2857
$ws.Range("L33").Value = 913.2857
$ws.Range("N33").Value = -1371.2857

# Row 41 (ALC) - hunk 1
$ws.Range("H41").Value = 1602.4166
$ws.Range("J41").Value = 1080
$ws.Range("L41").Value = 1080
$ws.Range("N41").Value = -1960

# Row 125 (ALC) - hunk 2
$ws.Range("H125").Value = 2798
$ws.Range("J125").Value = 3425.7144
$ws.Range("L125").Value = 30831.4296
$ws.Range("N125").Value = -35751.4296

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM) - hunk 3
$ws.Range("H2").Value = 3439.625
$ws.Range("I2").Value = 2672.3333
$ws.Range("K2").Value = 2672.3333
$ws.Range("M2").Value = -2559.3333

# Row 63 (ARM) - hunk 4
$ws.Range("H63").Value = 131109.1
$ws.Range("J63").Value = 204600
$ws.Range("L63").Value = 204600
$ws.Range("N63").Value = -205972

# Row 66 (ARM) - hunk 5
$ws.Range("H66").Value = 131109.1
$ws.Range("J66").Value = 204600
$ws.Range("L66").Value = 1023000
$ws.Range("N66").Value = -1029864

# Row 116 (ARM) - hunk 6
$ws.Range("H116").Value = 3439.625
$ws.Range("I116").Value = 2672.3333
$ws.Range("K116").Value = 2672.3333
$ws.Range("M116").Value = -378.3332999999998

# Row 122 (ARM) - hunk 7
$ws.Range("H122").Value = 5708.609
$ws.Range("I122").Value = 3984.6155
$ws.Range("J122").Value = 7949.8
$ws.Range("K122").Value = 11953.8465
$ws.Range("L122").Value = 23849.4
$ws.Range("M122").Value = -9503.8465
$ws.Range("N122").Value = -28749.4

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM) - hunk 8
$ws.Range("H3").Value = 3439.625
$ws.Range("I3").Value = 2672.3333
$ws.Range("K3").Value = 2672.3333
$ws.Range("M3").Value = -2558.3333

# Row 94 (BSM) - hunk 9
$ws.Range("H94").Value = 1411.1666
$ws.Range("I94").Value = 465.27777
$ws.Range("K94").Value = 465.27777
$ws.Range("M94").Value = -14.27776999999998

# Row 134 (BSM) - hunk 10
$ws.Range("H134").Value = 30305536
$ws.Range("I134").Value = 1868.5
$ws.Range("J134").Value = 66669936
$ws.Range("K134").Value = 5605.5
$ws.Range("L134").Value = 200009808
$ws.Range("M134").Value = -3070.5
$ws.Range("N134").Value = -200014878

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP) - hunk 11
$ws.Range("H31").Value = 3520.9167
$ws.Range("I31").Value = 1896.8334
$ws.Range("K31").Value = 1896.8334
$ws.Range("M31").Value = -1601.8334

# Row 34 (CRP) - hunk 12
$ws.Range("H34").Value = 3520.9167
$ws.Range("I34").Value = 1896.8334
$ws.Range("K34").Value = 1896.8334
$ws.Range("M34").Value = -1694.8334

# Row 62 (CRP) - hunk 13
$ws.Range("H62").Value = 490
$ws.Range("I62").Value = 490
$ws.Range("K62").Value = 490
$ws.Range("M62").Value = 134

# Row 65 (CRP) - hunk 14
$ws.Range("H65").Value = 490
$ws.Range("I65").Value = 490
$ws.Range("K65").Value = 2450
$ws.Range("M65").Value = 670

# Row 68 (CRP) - hunk 15
$ws.Range("H68").Value = 59995
$ws.Range("J68").Value = 59995
$ws.Range("L68").Value = 59995
$ws.Range("N68").Value = -61493

# Row 71 (CRP) - hunk 16
$ws.Range("H71").Value = 59995
$ws.Range("J71").Value = 59995
$ws.Range("L71").Value = 179985
$ws.Range("N71").Value = -187473

# Row 99 (CRP) - hunk 17
$ws.Range("H99").Value = 3217.5
$ws.Range("I99").Value = 3372.3635
$ws.Range("J99").Value = 1514
$ws.Range("K99").Value = 3372.3635
$ws.Range("L99").Value = 1514
$ws.Range("M99").Value = -1874.3635
$ws.Range("N99").Value = -4510

# Row 126 (CRP) - hunk 18
$ws.Range("H126").Value = 3217.5
$ws.Range("I126").Value = 3372.3635
$ws.Range("J126").Value = 1514
$ws.Range("K126").Value = 10117.0905
$ws.Range("L126").Value = 4542
$ws.Range("M126").Value = -7647.0905
$ws.Range("N126").Value = -9482

# Row 132 (CRP) - hunk 19
$ws.Range("H132").Value = 1583.4615
$ws.Range("I132").Value = 1465.4166
$ws.Range("K132").Value = 4396.2498
$ws.Range("M132").Value = -1866.2498

# Row 141 (CRP) - hunk 20
$ws.Range("H141").Value = 720403.9399999999
$ws.Range("J141").Value = 720403.9399999999
$ws.Range("L141").Value = 720403.9399999999
$ws.Range("N141").Value = -730763.9399999999

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (CUL) - hunk 21
$ws.Range("H4").Value = 84467256
$ws.Range("I4").Value = 81324650
$ws.Range("K4").Value = 243973950
$ws.Range("M4").Value = -243973838

# Row 96 (CUL) - hunk 22
$ws.Range("H96").Value = 2000
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 6000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -3941
$ws.Range("N96").ClearContents()

# Row 131 (CUL) - hunk 23
$ws.Range("H131").Value = 1475.5714
$ws.Range("I131").Value = 1070.4
$ws.Range("J131").Value = 1943.0769
$ws.Range("K131").Value = 3211.2
$ws.Range("L131").Value = 5829.2307
$ws.Range("M131").Value = 1828.8
$ws.Range("N131").Value = -15909.2307

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM) - hunk 24
$ws.Range("H70").Value = 4277.6
$ws.Range("I70").Value = 3249
$ws.Range("K70").Value = 3249
$ws.Range("M70").Value = -2979

# Row 73 (GSM) - hunk 25
$ws.Range("H73").Value = 4277.6
$ws.Range("I73").Value = 3249
$ws.Range("K73").Value = 3249
$ws.Range("M73").Value = -2313

# Row 97 (GSM) - hunk 26
$ws.Range("H97").Value = 2057.1428
$ws.Range("J97").Value = 3666.3333
$ws.Range("L97").Value = 3666.3333
$ws.Range("N97").Value = -4658.3333

# Row 107 (GSM) - hunk 27
$ws.Range("H107").Value = 379.89474
$ws.Range("I107").Value = 316.2
$ws.Range("J107").Value = 450.66666
$ws.Range("K107").Value = 316.2
$ws.Range("L107").Value = 450.66666
$ws.Range("M107").Value = 1603.8
$ws.Range("N107").Value = -4290.66666

# Row 122 (GSM) - hunk 28
$ws.Range("H122").Value = 1852.238
$ws.Range("I122").Value = 1557.2858
$ws.Range("K122").Value = 4671.857400000001
$ws.Range("M122").Value = -2221.857400000001

# Row 127 (GSM) - hunk 29
$ws.Range("H127").Value = 138775.33
$ws.Range("J127").Value = 138775.33
$ws.Range("L127").Value = 138775.33
$ws.Range("N127").Value = -148695.33

# Row 132 (GSM) - hunk 30
$ws.Range("H132").Value = 3014.5715
$ws.Range("I132").Value = 2850.3333
$ws.Range("K132").Value = 8550.999899999999
$ws.Range("M132").Value = -6020.999899999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW) - hunk 31
$ws.Range("H7").Value = 7531.8076
$ws.Range("I7").Value = 7210.0435
$ws.Range("K7").Value = 7210.0435
$ws.Range("M7").Value = -7098.0435

# Row 16 (LTW) - hunk 32
$ws.Range("H16").Value = 1851.4546
$ws.Range("I16").Value = 1851.4546
$ws.Range("K16").Value = 1851.4546
$ws.Range("M16").Value = -1681.4546

# Row 40 (LTW) - hunk 33
$ws.Range("H40").Value = 5217.727
$ws.Range("I40").Value = 4923.375
$ws.Range("K40").Value = 4923.375
$ws.Range("M40").Value = -4787.375

# Row 46 (LTW) - hunk 34
$ws.Range("H46").Value = 2178.4119
$ws.Range("I46").Value = 729.75
$ws.Range("K46").Value = 729.75
$ws.Range("M46").Value = -541.75

# Row 68 (LTW) - hunk 35
$ws.Range("H68").Value = 2609.6667
$ws.Range("I68").Value = 2609.6667
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2609.6667
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1860.6667
$ws.Range("N68").ClearContents()

# Row 71 (LTW) - hunk 36
$ws.Range("H71").Value = 2609.6667
$ws.Range("I71").Value = 2609.6667
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 13048.3335
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -9304.333500000001
$ws.Range("N71").ClearContents()

# Row 122 (LTW) - hunk 37
$ws.Range("H122").Value = 27719.2
$ws.Range("I122").Value = 22635.25
$ws.Range("K122").Value = 67905.75
$ws.Range("M122").Value = -65455.75

# Row 126 (LTW) - hunk 38
$ws.Range("H126").Value = 7531.8076
$ws.Range("I126").Value = 7210.0435
$ws.Range("K126").Value = 21630.1305
$ws.Range("M126").Value = -19160.1305

# Row 132 (LTW) - hunk 39
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 99 (WVR) - hunk 40
$ws.Range("H99").Value = 111000
$ws.Range("J99").Value = 111000
$ws.Range("L99").Value = 111000
$ws.Range("N99").Value = -116990

# Row 107 (WVR) - hunk 41
$ws.Range("H107").Value = 587
$ws.Range("I107").Value = 582
$ws.Range("J107").Value = 588.25
$ws.Range("K107").Value = 1746
$ws.Range("L107").Value = 1764.75
$ws.Range("M107").Value = 174
$ws.Range("N107").Value = -5604.75

# Row 126 (WVR) - hunk 42
$ws.Range("H126").Value = 21996
$ws.Range("I126").Value = 21996
$ws.Range("K126").Value = 65988
$ws.Range("M126").Value = -63518

# Row 132 (WVR) - hunk 43
$ws.Range("H132").Value = 2613.5625
$ws.Range("I132").Value = 2501.093
$ws.Range("K132").Value = 7503.279
$ws.Range("M132").Value = -4973.279
